$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in C2
$ws.Range("C2").Value = "Intro to scientific principles and data analysis"

# Shift the "Required Readings" column entries up by one, replacing the
# Clayton 2020 reference with a plain textbook citation, and trimming the
# Hu 2021 citation off the bottom entry.
$ws.Range("D2").Value = "LSWR Ch 2 and 3 <br>Light, Singer & Willet 1990, Ch. 2"
$ws.Range("D3").Value = "LSWR Ch 6 <br> [Clayton 2020](https://nautil.us/issue/92/frontiers/how-eugenics-shaped-statistics)"
$ws.Range("D4").Value = "LSWR Ch 11 and 12 <br> [Evans 2020](https://www.newstatesman.com/uncategorized/2020/07/ra-fisher-and-science-hatred)"
$ws.Range("D5").Value = "LSWR Ch 5 and 10 "

# Update the active selection on the sheet
$ws.Range("D3").Select()
